$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes
$ws.Columns.Item(1).ColumnWidth = 14.8325
$ws.Columns.Item(2).ColumnWidth = 15.66625

# Cell value changes
$ws.Range("A1").Value = -0.20264560746247184
$ws.Range("B1").Value = 0.20216905322129008
$ws.Range("A2").Value = -0.12841912878991213
$ws.Range("B2").Value = 0.1273926444852469
$ws.Range("A3").Value = -0.077677592124414829
$ws.Range("B3").Value = 0.077446960435734624
$ws.Range("A4").Value = -0.069446960467534069
$ws.Range("B4").Value = 0.069254152267941294
$ws.Range("A5").Value = -0.066254152284723311
$ws.Range("B5").Value = 0.065623420306144098
$ws.Range("A6").Value = -0.017351541342993215
$ws.Range("B6").Value = 0.017281649500565877
$ws.Range("A7").Value = -0.0072816495432976858
$ws.Range("B7").Value = 0.007276182123594932
$ws.Range("A8").Value = -0.009918556039143045
$ws.Range("B8").Value = 0.0098979729362658198
$ws.Range("A9").Value = -0.007897972953295973
$ws.Range("B9").Value = 0.0078879613018592032
$ws.Range("A10").Value = -0.0058879613191251678
$ws.Range("B10").Value = 0.005888117302161433
$ws.Range("A11").Value = -0.0028881173227111034
$ws.Range("B11").Value = 0.0028882057002412864
$ws.Range("A12").Value = 0.00061179427759983085
$ws.Range("B12").Value = -0.00061339808560312292
$ws.Range("A13").Value = 0.0041133980637315659
$ws.Range("B13").Value = -0.0041192138301795467
$ws.Range("A14").Value = 0.012119213793807759
$ws.Range("B14").Value = -0.012141734652164615
$ws.Range("A15").Value = -0.0080516403726456076
$ws.Range("B15").Value = 0.0080335433247222809
$ws.Range("A16").Value = -0.0060335433412688211
$ws.Range("B16").Value = 0.0060034935430643976
$ws.Range("A17").Value = -0.0040034935599511101
$ws.Range("B17").Value = 0.0039999999764495087
$ws.Range("A18").Value = -0.016105603993548101
$ws.Range("B18").Value = 0.016091948126934597
$ws.Range("A19").Value = -0.012091948140999786
$ws.Range("B19").Value = 0.012016991881657724
$ws.Range("A20").Value = -0.0080169918967101239
$ws.Range("B20").Value = 0.0080056983753049593
$ws.Range("A21").Value = -0.004005698390534107
$ws.Range("B21").Value = 0.0039999999846420664
$ws.Range("A22").Value = -0.045715052378083598
$ws.Range("B22").Value = 0.045500695461669238
$ws.Range("A23").Value = -0.040500695483331128
$ws.Range("B23").Value = 0.040099407218699312
$ws.Range("A24").Value = -0.020099407291542803
$ws.Range("B24").Value = 0.019999999926148426
$ws.Range("A25").Value = -0.046271878979732506
$ws.Range("B25").Value = 0.046239691393108018
$ws.Range("A26").Value = -0.043739691411744275
$ws.Range("B26").Value = 0.043701652074624064
$ws.Range("A27").Value = -0.041201652093859398
$ws.Range("B27").Value = 0.04099625772540838
$ws.Range("A28").Value = -0.038996257745589347
$ws.Range("B28").Value = 0.038870630823025998
$ws.Range("A29").Value = -0.031870630861465798
$ws.Range("B29").Value = 0.03184493170479552
$ws.Range("A30").Value = 0.028155068082261536
$ws.Range("B30").Value = -0.028206763507530486
$ws.Range("A31").Value = -0.014023234869798173
$ws.Range("B31").Value = 0.014001406425977692
$ws.Range("A32").Value = -0.0040014064737210475
$ws.Range("B32").Value = 0.0039999999718425272
